$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 343, shifting the existing rows (343:366) down to (344:367)
$ws.Rows("343:343").Insert()

# Populate the newly inserted row 343 with the new price-report entry.
# Columns A,B,C,E-L carry forward the same "Femacal de La Calera / Coquimbo / Arándano (blue)"
# metadata as the row that used to occupy this slot; D and M-T hold the new observation.
$ws.Range("A343").Value = 3
$ws.Range("B343").Value = "Femacal de La Calera"
$ws.Range("C343").Value = "Coquimbo"
$ws.Range("D343").Value = 45166
$ws.Range("E343").Value = 5
$ws.Range("F343").Value = "Fruta"
$ws.Range("G343").Value = 100101
$ws.Range("H343").Value = "Berries"
$ws.Range("I343").Value = 100101001
$ws.Range("J343").Value = "Arándano (blue)"
$ws.Range("K343").Value = "Sin especificar"
$ws.Range("L343").Value = "Primera"
$ws.Range("M343").Value = 45
$ws.Range("N343").Value = 13000
$ws.Range("O343").Value = 13000
$ws.Range("P343").Value = 13000
$ws.Range("Q343").Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Range("R343").Value = "Provincia de Quillota"
$ws.Range("S343").Value = 8667
$ws.Range("T343").Value = 1.5
